$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$years = 2005..2015
for ($i = 0; $i -lt $years.Count; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value2 = $years[$i]
}

$ws1.Activate()
$ws1.Range("A1:C12").Select()

$ps = $ws1.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
